# Travis County 2018 bg SVI workbook update
# Re-orders the factor-analysis "Significant Components" / "Loading Factors" /
# "All Refactor Variances" / "Final Variances" / "Included and Excluded"
# tables to reflect a re-run of the underlying analysis (same variables,
# slightly different loadings/variances and row ordering).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Significant Components"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Significant Components")

$ws1.Range("C2").Value = "['QSERV' 'QEXTRCT' 'QESL' 'QHISPC' 'QEDLESHI' 'PPUNIT' 'QNOHLTH' 'QFHH'`n 'PERCAP']"
$ws1.Range("C3").Value = "['QRICH' 'PERCAP' 'MDHSEVAL']"
$ws1.Range("C5").Value = "['QAGEDEP' 'MEDAGE' 'QSSBEN']"
$ws1.Range("C6").Value = "['QAGEDEP' 'QFEMLBR' 'QFEMALE']"

# The C2 text spans two lines; writing it nudges the row to a custom height
# in this runtime. Auto-fit it back so the row keeps its original (default)
# height, matching the source workbook.
$ws1.Rows.Item(2).AutoFit()

# ---------------------------------------------------------------------------
# Sheet "Loading Factors"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Loading Factors")

$ws2.Cells.Item(2,1).Value = "QSERV"
$ws2.Cells.Item(2,2).Value = 0.5274953379585156
$ws2.Cells.Item(2,3).Value = 0.3608946325620632
$ws2.Cells.Item(2,4).Value = 0.3130727392093795
$ws2.Cells.Item(2,5).Value = -0.1262082763595171
$ws2.Cells.Item(2,6).Value = -0.10392483115886

$ws2.Cells.Item(3,1).Value = "QEXTRCT"
$ws2.Cells.Item(3,2).Value = 0.786144157245031
$ws2.Cells.Item(3,3).Value = 0.1131783495998638
$ws2.Cells.Item(3,4).Value = 0.08480976424482217
$ws2.Cells.Item(3,5).Value = -0.02097332736471048
$ws2.Cells.Item(3,6).Value = -0.2028241015362909

$ws2.Cells.Item(4,1).Value = "QESL"
$ws2.Cells.Item(4,2).Value = 0.7994959196170861
$ws2.Cells.Item(4,3).Value = 0.1373791394723209
$ws2.Cells.Item(4,4).Value = 0.211719776206982
$ws2.Cells.Item(4,5).Value = -0.0326793945777714
$ws2.Cells.Item(4,6).Value = -0.2000727787374849

$ws2.Cells.Item(5,1).Value = "QHISPC"
$ws2.Cells.Item(5,2).Value = 0.8230600188201528
$ws2.Cells.Item(5,3).Value = 0.3465122419334931
$ws2.Cells.Item(5,4).Value = 0.1359332285180322
$ws2.Cells.Item(5,5).Value = -0.09949764747652885
$ws2.Cells.Item(5,6).Value = -0.1288344851589302

$ws2.Cells.Item(6,1).Value = "QEDLESHI"
$ws2.Cells.Item(6,2).Value = 0.8778539376176244
$ws2.Cells.Item(6,3).Value = 0.196185898045183
$ws2.Cells.Item(6,4).Value = 0.1910468545985881
$ws2.Cells.Item(6,5).Value = -0.01699417353382815
$ws2.Cells.Item(6,6).Value = -0.1019276631332453

$ws2.Cells.Item(7,1).Value = "PPUNIT"
$ws2.Cells.Item(7,2).Value = 0.7145302556558201
$ws2.Cells.Item(7,3).Value = 0.04770843604974272
$ws2.Cells.Item(7,4).Value = -0.3603295215360905
$ws2.Cells.Item(7,5).Value = -0.09691264601002293
$ws2.Cells.Item(7,6).Value = 0.07246473358154387

$ws2.Cells.Item(8,1).Value = "QNOHLTH"
$ws2.Cells.Item(8,2).Value = 0.6741858254799591
$ws2.Cells.Item(8,3).Value = 0.3959044364447006
$ws2.Cells.Item(8,4).Value = 0.324282530897769
$ws2.Cells.Item(8,5).Value = -0.08306123956933575
$ws2.Cells.Item(8,6).Value = -0.1222630309176319

# Row 9 (QFHH) keeps its label, only the values get a tiny refresh.
$ws2.Cells.Item(9,2).Value = 0.5602157023173483
$ws2.Cells.Item(9,3).Value = 0.2825588698036734
$ws2.Cells.Item(9,4).Value = 0.04381734165852227
$ws2.Cells.Item(9,5).Value = -0.06322069616271428
$ws2.Cells.Item(9,6).Value = 0.2438533305152224

$ws2.Cells.Item(10,1).Value = "QRICH"
$ws2.Cells.Item(10,2).Value = 0.216768309646086
$ws2.Cells.Item(10,3).Value = 0.8465873572509584
$ws2.Cells.Item(10,4).Value = 0.3549485723718843
$ws2.Cells.Item(10,5).Value = -0.1549872521244396
$ws2.Cells.Item(10,6).Value = -0.04401484876180284

$ws2.Cells.Item(11,1).Value = "PERCAP"
$ws2.Cells.Item(11,2).Value = 0.474525869952907
$ws2.Cells.Item(11,3).Value = 0.7338494751353548
$ws2.Cells.Item(11,4).Value = 0.247403648082836
$ws2.Cells.Item(11,5).Value = -0.2041762660230629
$ws2.Cells.Item(11,6).Value = 0.03151731727370427

$ws2.Cells.Item(12,1).Value = "MDHSEVAL"
$ws2.Cells.Item(12,2).Value = 0.3004626931611233
$ws2.Cells.Item(12,3).Value = 0.7871221240119082
$ws2.Cells.Item(12,4).Value = -0.06641551582182338
$ws2.Cells.Item(12,5).Value = -0.0496803877291036
$ws2.Cells.Item(12,6).Value = 0.02437791517645226

# Row 13 (QRENTER) keeps its label, only the values get a tiny refresh.
$ws2.Cells.Item(13,2).Value = -0.005926326639922076
$ws2.Cells.Item(13,3).Value = 0.1780075773328002
$ws2.Cells.Item(13,4).Value = 0.7703798488407927
$ws2.Cells.Item(13,5).Value = -0.4138178469644066
$ws2.Cells.Item(13,6).Value = -0.1132850040875337

# Row 14 (QNOAUTO) keeps its label, only the values get a tiny refresh.
$ws2.Cells.Item(14,2).Value = 0.1139529751834322
$ws2.Cells.Item(14,3).Value = 0.06106384712282543
$ws2.Cells.Item(14,4).Value = 0.6847406410186511
$ws2.Cells.Item(14,5).Value = -0.05243788577041029
$ws2.Cells.Item(14,6).Value = 0.007215023145616352

# Row 15 (QPOVTY) keeps its label, only the values get a tiny refresh.
$ws2.Cells.Item(15,2).Value = 0.3023780544425386
$ws2.Cells.Item(15,3).Value = 0.1695709324950138
$ws2.Cells.Item(15,4).Value = 0.5567092088776491
$ws2.Cells.Item(15,5).Value = -0.2796814514076003
$ws2.Cells.Item(15,6).Value = 0.09476063937229715

$ws2.Cells.Item(16,1).Value = "QAGEDEP"
$ws2.Cells.Item(16,2).Value = -0.01897725464058234
$ws2.Cells.Item(16,3).Value = -0.1465327652254485
$ws2.Cells.Item(16,4).Value = -0.1088299664760572
$ws2.Cells.Item(16,5).Value = 0.6889723574127111
$ws2.Cells.Item(16,6).Value = 0.6048116232004375

$ws2.Cells.Item(17,1).Value = "MEDAGE"
$ws2.Cells.Item(17,2).Value = -0.3203145699989389
$ws2.Cells.Item(17,3).Value = -0.2374924602177076
$ws2.Cells.Item(17,4).Value = -0.3581218318328633
$ws2.Cells.Item(17,5).Value = 0.7551459380736495
$ws2.Cells.Item(17,6).Value = -0.03604536740633013

# Row 18 (QSSBEN) keeps its label, only the values get a tiny refresh.
$ws2.Cells.Item(18,2).Value = 0.006717323940582034
$ws2.Cells.Item(18,3).Value = -0.04185582841690016
$ws2.Cells.Item(18,4).Value = -0.1566751810311308
$ws2.Cells.Item(18,5).Value = 0.8301134378939941
$ws2.Cells.Item(18,6).Value = 0.08972562158399182

$ws2.Cells.Item(19,1).Value = "QFEMLBR"
$ws2.Cells.Item(19,2).Value = -0.2650015962054501
$ws2.Cells.Item(19,3).Value = 0.084970300904067
$ws2.Cells.Item(19,4).Value = 0.005596367903975159
$ws2.Cells.Item(19,5).Value = -0.04604931457693156
$ws2.Cells.Item(19,6).Value = 0.7314609960898517

$ws2.Cells.Item(20,1).Value = "QFEMALE"
$ws2.Cells.Item(20,2).Value = -0.0188953065370911
$ws2.Cells.Item(20,3).Value = -0.05607244907176379
$ws2.Cells.Item(20,4).Value = 0.00262935355599888
$ws2.Cells.Item(20,5).Value = 0.1873923085423901
$ws2.Cells.Item(20,6).Value = 0.8727849509922407

# ---------------------------------------------------------------------------
# Sheet "All Refactor Variances" - only the "F2" block (columns I:R) refreshes
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("All Refactor Variances")

$ws3.Cells.Item(2,9).Value  = 4.751932265115403
$ws3.Cells.Item(2,10).Value = 3.118803383093232
$ws3.Cells.Item(2,11).Value = 2.379793064865748
$ws3.Cells.Item(2,12).Value = 2.14111305618065
$ws3.Cells.Item(2,13).Value = 1.88606688715619
$ws3.Cells.Item(2,14).Value = 4.902928252043285
$ws3.Cells.Item(2,15).Value = 2.58830961334677
$ws3.Cells.Item(2,16).Value = 2.17049435470839
$ws3.Cells.Item(2,17).Value = 2.139559849957482
$ws3.Cells.Item(2,18).Value = 1.895944382860769

$ws3.Cells.Item(3,9).Value  = 0.2262824888150192
$ws3.Cells.Item(3,10).Value = 0.1485144468139634
$ws3.Cells.Item(3,11).Value = 0.1133234792793213
$ws3.Cells.Item(3,12).Value = 0.1019577645800309
$ws3.Cells.Item(3,13).Value = 0.08981270891219954
$ws3.Cells.Item(3,14).Value = 0.2580488553706992
$ws3.Cells.Item(3,15).Value = 0.1362268217550932
$ws3.Cells.Item(3,16).Value = 0.1142365449846521
$ws3.Cells.Item(3,17).Value = 0.112608413155657
$ws3.Cells.Item(3,18).Value = 0.09978654646635626

$ws3.Cells.Item(4,9).Value  = 0.2262824888150192
$ws3.Cells.Item(4,10).Value = 0.3747969356289826
$ws3.Cells.Item(4,11).Value = 0.4881204149083039
$ws3.Cells.Item(4,12).Value = 0.5900781794883349
$ws3.Cells.Item(4,13).Value = 0.6798908884005344
$ws3.Cells.Item(4,14).Value = 0.2580488553706992
$ws3.Cells.Item(4,15).Value = 0.3942756771257924
$ws3.Cells.Item(4,16).Value = 0.5085122221104444
$ws3.Cells.Item(4,17).Value = 0.6211206352661014
$ws3.Cells.Item(4,18).Value = 0.7209071817324577

$ws3.Cells.Item(5,9).Value  = 0.3328217698980437
$ws3.Cells.Item(5,10).Value = 0.218438648535721
$ws3.Cells.Item(5,11).Value = 0.1666789204160664
$ws3.Cells.Item(5,12).Value = 0.1499619517182969
$ws3.Cells.Item(5,13).Value = 0.1320987094318721
$ws3.Cells.Item(5,14).Value = 0.3579501798699878
$ws3.Cells.Item(5,15).Value = 0.1889658269566935
$ws3.Cells.Item(5,16).Value = 0.1584622096704919
$ws3.Cells.Item(5,17).Value = 0.1562037610515137
$ws3.Cells.Item(5,18).Value = 0.1384180224513132

# ---------------------------------------------------------------------------
# Sheet "Final Variances" (mirrors the "F2" block above)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Final Variances")

$ws4.Cells.Item(2,2).Value = 4.902928252043285
$ws4.Cells.Item(2,3).Value = 2.58830961334677
$ws4.Cells.Item(2,4).Value = 2.17049435470839
$ws4.Cells.Item(2,5).Value = 2.139559849957482
$ws4.Cells.Item(2,6).Value = 1.895944382860769

$ws4.Cells.Item(3,2).Value = 0.2580488553706992
$ws4.Cells.Item(3,3).Value = 0.1362268217550932
$ws4.Cells.Item(3,4).Value = 0.1142365449846521
$ws4.Cells.Item(3,5).Value = 0.112608413155657
$ws4.Cells.Item(3,6).Value = 0.09978654646635626

$ws4.Cells.Item(4,2).Value = 0.2580488553706992
$ws4.Cells.Item(4,3).Value = 0.3942756771257924
$ws4.Cells.Item(4,4).Value = 0.5085122221104444
$ws4.Cells.Item(4,5).Value = 0.6211206352661014
$ws4.Cells.Item(4,6).Value = 0.7209071817324577

$ws4.Cells.Item(5,2).Value = 0.3579501798699878
$ws4.Cells.Item(5,3).Value = 0.1889658269566935
$ws4.Cells.Item(5,4).Value = 0.1584622096704919
$ws4.Cells.Item(5,5).Value = 0.1562037610515137
$ws4.Cells.Item(5,6).Value = 0.1384180224513132

# ---------------------------------------------------------------------------
# Sheet "Included and Excluded"
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Included and Excluded")

$ws5.Range("B2").Value = "[['QSERV', 'QEXTRCT', 'QESL', 'QHISPC', 'QEDLESHI', 'PPUNIT', 'QNOHLTH', 'QFHH', 'PERCAP', 'QRICH', 'MDHSEVAL', 'QRENTER', 'QNOAUTO', 'QPOVTY', 'QAGEDEP', 'MEDAGE', 'QSSBEN', 'QFEMLBR', 'QFEMALE']]"
